$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.536.34"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.902.61"
$ws.Range("E3").Value = "  +3.09%  "
$ws.Range("E4").Value = "  +0.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.14"
$ws.Range("E5").Value = "  +5.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.631"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.16"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("E9").Value = "  +2.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0703"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0997"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.178.44"
$ws.Range("E12").Value = "  +3.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.38"
$ws.Range("E13").Value = "  +7.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.912.41"
$ws.Range("E14").Value = "  +3.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.689"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.83"
$ws.Range("E16").Value = "  +2.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.519.01"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "71.82"
$ws.Range("E18").Value = "  +2.49%  "
$ws.Range("D19").Value = "0.0₃0811"
$ws.Range("E19").Value = "  +2.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "243.37"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.87"
$ws.Range("E22").Value = "  +2.07%  "
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.26"
$ws.Range("E25").Value = "  +31.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.80"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.54"
$ws.Range("E27").Value = "  +8.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.95"
$ws.Range("E28").Value = "  +2.46%  "
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.972"
$ws.Range("E30").Value = "  +27.33%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0565"
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.09"
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("E34").Value = "  +4.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.75"
$ws.Range("E35").Value = "  +6.46%  "
$ws.Range("E36").Value = "  +2.74%  "
$ws.Range("E37").Value = "  +4.16%  "
$ws.Range("E38").Value = "  +3.15%  "
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "90.97"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.05"
$ws.Range("E41").Value = "  +48.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.350.19"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.50"
$ws.Range("E43").Value = "  +6.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0591"
$ws.Range("E44").Value = "  +11.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.33"
$ws.Range("E45").Value = "  +1.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.63"
$ws.Range("E46").Value = "  +8.92%  "
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.64"
$ws.Range("E49").Value = "  +4.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.087.37"
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("E51").Value = "  +2.43%  "
